# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Updates the "K" column (column G) values on the active worksheet to the
# newly-computed strikeout counts for each row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (strikeout) value
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 3
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 2
    18 = 3
    19 = 4
    20 = 1
    21 = 0
    22 = 4
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
